# Applies the "Atualização de bases das ligas, do dia: 04-04-2024 às 23:22" update
# to the "Australia ALeague" sheet.
#
# Two kinds of changes:
#  1) Three pairs of data rows got swapped in full (every column from B to AC,
#     i.e. everything except the running index in column A): 104<->105,
#     112<->113, 124<->125.
#  2) A handful of odds cells (columns N,O,P,Q,R,S,T,U,V) were corrected in
#     rows 139-144, with no row movement.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA`:AC$rowA")
    $rangeB = $ws.Range("B$rowB`:AC$rowB")

    $valsA = $rangeA.Value2
    $valsB = $rangeB.Value2

    $rangeA.Value = $valsB
    $rangeB.Value = $valsA
}

Swap-Rows 104 105
Swap-Rows 112 113
Swap-Rows 124 125

# In-place odds corrections for rows 139-144 (no row movement here).
$ws.Range("N139").Value = 2.1
$ws.Range("P139").Value = 3.1
$ws.Range("Q139").Value = -0.25
$ws.Range("R139").Value = 1.85
$ws.Range("S139").Value = 2.05

$ws.Range("P140").Value = 3.3

$ws.Range("N141").Value = 1.75
$ws.Range("O141").Value = 4
$ws.Range("P141").Value = 4.5
$ws.Range("R141").Value = 2.01
$ws.Range("S141").Value = 1.89
$ws.Range("U141").Value = 2.025
$ws.Range("V141").Value = 1.825

$ws.Range("N142").Value = 2.25
$ws.Range("O142").Value = 3.6
$ws.Range("P142").Value = 3.1
$ws.Range("R142").Value = 1.98
$ws.Range("S142").Value = 1.92
$ws.Range("T142").Value = 3
$ws.Range("U142").Value = 1.975
$ws.Range("V142").Value = 1.875

$ws.Range("N143").Value = 4.2
$ws.Range("O143").Value = 4
$ws.Range("Q143").Value = 0.75
$ws.Range("R143").Value = 1.89
$ws.Range("S143").Value = 2.01
$ws.Range("U143").Value = 1.9
$ws.Range("V143").Value = 1.95

$ws.Range("N144").Value = 2.8
$ws.Range("P144").Value = 2.45
$ws.Range("R144").Value = 2.07
$ws.Range("S144").Value = 1.83
$ws.Range("U144").Value = 1.925
$ws.Range("V144").Value = 1.925

Write-Output "Applied Australia ALeague update."
